# Automatic update of files.
# Update the "Förändrad" (Changed) date column (C) for rows 2-42 from
# 2025-02-24 (45712) to 2025-02-25 (45713) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C42").Value = 45713
